$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.91
$ws.Range("Q2").Value = 1.5
$ws.Range("R2").Value = 2.63
$ws.Range("AM2").Value = 26
